$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5708487629890442
$ws.Range("B1").Value = 0.8488417863845825
$ws.Range("C1").Value = 4.541937828063965
$ws.Range("D1").Value = 2.075095653533936
$ws.Range("E1").Value = 0.7904588580131531
